## Results of ConQat.xlsx - "Result updated for Jabref, RQ1."
##
## RQ1 sheet: the MonoOSC row is removed from the small results table
## (rows 3-8/9) and the Jabref row is filled in with real CF/Revisions
## data (it previously only had a "Revisions" count with no CF data).
## Rows below the table (23-30) are untouched, so this is done with
## direct cell edits rather than a sheet/row delete (which would shift
## everything below it).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RQ1")
$ws4 = $wb.Worksheets.Item("RQ4")

# --- Row 5: was MonoOSC/1110, becomes Freecol/10354 --------------------
$ws1.Range("A5").Value = "Freecol"
$ws1.Range("B5").Value = 10354
$ws1.Range("E5").Formula = "=C5/D5"
$ws1.Range("H5").Formula = "=F5/G5"

# --- Row 6: was Freecol/10354, becomes Carol/2886 -----------------------
$ws1.Range("A6").Value = "Carol"
$ws1.Range("B6").Value = 2886
$ws1.Range("E6").Formula = "=C6/D6"
$ws1.Range("H6").Formula = "=F6/G6"

# --- Row 7: was Carol/2886 (no CF data), becomes Jabref/2798 with data --
$ws1.Range("A7").Value = "Jabref"
$ws1.Range("B7").Value = 2798
$ws1.Range("C7").Value = 165
$ws1.Range("D7").Value = 148
$ws1.Range("E7").Formula = "=C7/D7"
$ws1.Range("F7").Value = 1804
$ws1.Range("G7").Value = 148
$ws1.Range("H7").Formula = "=F7/G7"

# --- Row 8: was Jabref/2798 (no CF data), becomes the Total row --------
$ws1.Range("A8").Value = "Total"
$ws1.Range("B8").Value = $null
$ws1.Range("C8").Formula = "=SUM(C3:C7)"
$ws1.Range("D8").Formula = "=SUM(D3:D7)"
$ws1.Range("E8").Formula = "=C8/D8"
$ws1.Range("F8").Formula = "=SUM(F3:F7)"
$ws1.Range("G8").Formula = "=SUM(G3:G7)"
$ws1.Range("H8").Formula = "=F8/G8"

# --- Old row 9 (the previous Total row) is now gone completely; clear --
# --- it fully (contents + formatting) so no empty row 9 remains.       --
$ws1.Range("A9:H9").Clear()

# --- Selection / active sheet: RQ1 becomes the active tab, selection ---
# --- on RQ1 moves to E7; RQ4 is no longer the selected tab.            --
[void]$ws4.Activate()
[void]$ws4.Range("F3").Select()
[void]$ws1.Activate()
[void]$ws1.Range("E7").Select()
